# Actualización automática del tracker
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resultado / profit updates for rows that were pending ("Fallo" / "Acierto")
$updates = @{
    293 = @("Fallo", -1)
    294 = @("Fallo", -1)
    333 = @("Fallo", -1)
    334 = @("Fallo", -1)
    335 = @("Fallo", -1)
    336 = @("Fallo", -1)
    338 = @("Fallo", -1)
    343 = @("Fallo", -1)
    344 = @("Acierto", 2.5)
    345 = @("Fallo", -1)
    346 = @("Acierto", 2.4)
    355 = @("Fallo", -1)
    356 = @("Fallo", -1)
    363 = @("Acierto", 1.1)
    367 = @("Fallo", -1)
    369 = @("Fallo", -1)
    370 = @("Fallo", -1)
    371 = @("Fallo", -1)
    373 = @("Acierto", 2)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("G$row").Value = $vals[0]
    $ws.Range("H$row").Value = $vals[1]
}

# Append new row 376 with the latest scraped result
$ws.Range("A376").Value = 14427995

# "fecha" column stores dates as plain text (e.g. "2025-08-21"), not native
# Excel dates -- force text formatting so COM doesn't auto-convert the
# string into a date serial number.
$ws.Range("B376").NumberFormat = "@"
$ws.Range("B376").Value = "2025-08-21"
$ws.Range("B376").Style = "Normal"

$ws.Range("C376").Value = "Liudmila Samsonova"
$ws.Range("D376").Value = "Sorana Cirstea"
$ws.Range("E376").Value = "Gana Liudmila Samsonova"
$ws.Range("F376").Value = 1.57
$ws.Range("G376").Value = "Fallo"
$ws.Range("H376").Value = -1
